# --------------------------------------------------------------------------
# PlayerPerformance_4601.xlsx edit script
#
# Starting point (before.xlsx):
#   Sheet 1: "ODI Batting"   (MATCH_NUMBER, INNING_NUMBER, MATCH_DATE, MATCH_CARD_LINK, ...)
#   Sheet 2: "ODI Bowling"   (MATCH_NUMBER, MATCH_CARD_LINK, MATCH_INNING, ...)
#
# Target (after edit):
#   Sheet 1: "Player Info"        (NEW) - ID, NAME, BATTING_HAND, BOWL_STYLE
#   Sheet 2: "ODI Batting"        (existing, MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code,
#                                  blank INNING_NUMBER cells removed)
#   Sheet 3: "ODI Bowling"        (existing, MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code)
#   Sheet 4: "ODI Batting Extra"  (NEW) - MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
#                                  PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the same bold/centered/bordered header look used by the
# existing header rows (style index 1 in the original file) to a range.
# ---------------------------------------------------------------------
function Format-HeaderRange($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------
# STEP 1: rebuild the sheet order / names
#   insert "Player Info" before "ODI Batting"
#   insert "ODI Batting Extra" after "ODI Bowling"
# ---------------------------------------------------------------------
$wsBattingOrig = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($wsBattingOrig)
$playerInfo.Name = "Player Info"

$wsBowlingOrig = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add($null, $wsBowlingOrig)
$battingExtra.Name = "ODI Batting Extra"

Write-Host "Sheet order:"
foreach ($s in $wb.Worksheets) { Write-Host " - " $s.Name }

# ---------------------------------------------------------------------
# STEP 2: populate "Player Info"
# ---------------------------------------------------------------------
$wsPlayer = $wb.Worksheets.Item("Player Info")

$playerHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerHeaders.Length; $c++) {
    $wsPlayer.Cells.Item(1, $c).Value = $playerHeaders[$c - 1]
}
Format-HeaderRange $wsPlayer.Range("A1:D1")

# keep the player id as text (matches source data convention of numeric-looking codes as text)
$wsPlayer.Range("A2").NumberFormat = "@"
$wsPlayer.Cells.Item(2, 1).Value = "4601"
$wsPlayer.Cells.Item(2, 2).Value = "Keshav A Maharaj"
$wsPlayer.Cells.Item(2, 3).Value = "Right Handed"
$wsPlayer.Cells.Item(2, 4).Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------
# STEP 3: fix up "ODI Batting"
#   D1 header: MATCH_CARD_LINK -> MATCH_CODE
#   D2:D28   : full scorecard URL -> bare numeric match code (kept as text)
#   B column : rows with a blank INNING_NUMBER lose the cell entirely
# ---------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsBatting.Range("D1").Value = "MATCH_CODE"

$wsBatting.Range("D2:D28").NumberFormat = "@"
for ($r = 2; $r -le 28; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url) {
        $code = $url -replace '^.*MatchCode=', ''
        $cell.Value = $code
    }
}

$blankInningRows = @(2, 3, 7, 8, 10, 12, 14, 17, 18, 21, 23, 25, 26, 28)
foreach ($r in $blankInningRows) {
    $wsBatting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------
# STEP 4: fix up "ODI Bowling"
#   B1 header : MATCH_CARD_LINK -> MATCH_CODE
#   B2:B27    : full scorecard URL -> bare numeric match code (kept as text)
# ---------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

$wsBowling.Range("B1").Value = "MATCH_CODE"

$wsBowling.Range("B2:B27").NumberFormat = "@"
for ($r = 2; $r -le 27; $r++) {
    $cell = $wsBowling.Cells.Item($r, 2)
    $url = $cell.Value2
    if ($url) {
        $code = $url -replace '^.*MatchCode=', ''
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------
# STEP 5: build the new "ODI Batting Extra" sheet
#   MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
#   BATTING_POSITION is numeric; everything else (incl. blanks) is text.
# ---------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $wsExtra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
Format-HeaderRange $wsExtra.Range("A1:F1")

# keep match codes as text
$wsExtra.Range("A2:A21").NumberFormat = "@"
# the rest of the non-numeric columns (blank cells mix with text) stay as text too
$wsExtra.Range("C2:F21").NumberFormat = "@"

$extraRows = @(
    @("4460", 8,    "3", "0", "6.51%",  "NO"),
    @("4474", $null, "",  "",  "",       "NO"),
    @("4475", 8,    "1", "0", "6.88%",  "NO"),
    @("4478", $null, "",  "",  "",       "NO"),
    @("4487", 9,    "0", "0", "",       "NO"),
    @("4488", 10,   "",  "",  "",       "NO"),
    @("4491", 10,   "2", "0", "12.00%", "NO"),
    @("4517", $null, "",  "",  "",       "NO"),
    @("4524", $null, "",  "",  "",       "NO"),
    @("4526", 8,    "",  "",  "",       "NO"),
    @("4529", 9,    "1", "0", "2.09%",  "NO"),
    @("4550", 10,   "3", "1", "8.33%",  "NO"),
    @("4557", 9,    "",  "",  "",       "NO"),
    @("4559", 8,    "4", "0", "18.18%", "NO"),
    @("4619", $null, "",  "",  "",       "NO"),
    @("4620", 8,    "0", "0", "1.20%",  "NO"),
    @("4622", $null, "",  "",  "",       "NO"),
    @("4656", $null, "",  "",  "",       "NO"),
    @("4657", 8,    "0", "0", "1.80%",  "NO"),
    @("4699", 9,    "",  "",  "",       "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $wsExtra.Cells.Item($r, 1).Value = $row[0]

    if ($null -eq $row[1]) {
        $wsExtra.Cells.Item($r, 2).ClearContents()
    } else {
        $wsExtra.Cells.Item($r, 2).Value = $row[1]
    }

    for ($c = 3; $c -le 6; $c++) {
        $val = $row[$c - 1]
        if ($val -eq "") {
            $wsExtra.Cells.Item($r, $c).Value = ""
        } else {
            $wsExtra.Cells.Item($r, $c).Value = $val
        }
    }

    $r = $r + 1
}

# ---------------------------------------------------------------------
# STEP 6: restore the original active-tab (tab 0 / "Player Info", the
# new first sheet) since adding sheets shifts Excel's active selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Player Info").Activate()
